$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet (LOGT1 -> logt1)
$ws.Name = "logt1"

# 2. Update the assembly date at F1 (41898 -> 41902)
$ws.Range("F1").Value = 41902

# 3. Remove the old totals row (row 11) - its formula is rebuilt at the new
#    location (row 5) further down, right above the header.
$ws.Rows.Item(11).Delete()

# 4. Insert a blank row above the old header row (row 5), pushing the header
#    and all log rows down by one. This new row 5 will hold the total
#    formula.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).RowHeight = 13.75

# 5. Insert two new blank rows between the two "github" log entries (now at
#    rows 7 and 8) and the two "ruby" log entries (now at rows 9 and 10), to
#    host two brand-new log entries.
$ws.Range("A9:A10").EntireRow.Insert()
$ws.Rows.Item(9).RowHeight = 14.15
$ws.Rows.Item(10).RowHeight = 14.15

# 6. Fill in the two new log entries.
$ws.Range("A9").Value = 41902
$ws.Range("B9").Value = 91/144
$ws.Range("C9").Value = 197/288
$ws.Range("D9").Value = 14
$ws.Range("E9").Formula = "=((HOUR(C9)-HOUR(B9))*60)+(MINUTE(C9)-MINUTE(B9))-D9"
$ws.Range("F9").Value = 2
$ws.Range("H9").Value = "Realizar el lanzamiento del ciclo #1 de TSPi."

$ws.Range("A10").Value = 41902
$ws.Range("B10").Value = 11/16
$ws.Range("C10").Value = 209/288
$ws.Range("D10").Value = 10
$ws.Range("E10").Formula = "=((HOUR(C10)-HOUR(B10))*60)+(MINUTE(C10)-MINUTE(B10))-D10"
$ws.Range("F10").Value = 3
$ws.Range("H10").Value = "Definir la estrategía de desarrollo del ciclo #1 de TSPi."

# 7. Replace the "-" placeholders in F/G for the first two log rows (7, 8)
#    and the last two log rows (11, 12) with the new sequence numbers,
#    clearing the G column (Assembly) for all the log rows.
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = ""
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = ""
$ws.Range("F11").Value = 11
$ws.Range("G11").Value = ""
$ws.Range("F12").Value = 11
$ws.Range("G12").Value = ""

# 8. Right-align the Phase/Task numbers in column F for every log row.
$ws.Range("F7:F12").HorizontalAlignment = -4152

# 9. Add the new total formula at E5.
$ws.Range("E5").Formula = "=SUM(E7:E12)/60"

# 10. Restore the selected cell to E6 (the new header row).
$ws.Range("E6").Select()
